# ADC calcs rework: re-configured ADC to use internal Vref, added filtering
# for ADC data, added hysteresis for ADC thresholds.
#
# - A1 (Vref-derived ratio): 0.17 -> 0.55000000000000004
# - A3 (numerator constant): 11 -> 10.65
# - A4 (filter weight): 10 -> 1
# - A6 (hysteresis divisor): 5.0999999999999996 -> 1.1000000000000001
# - Column D (the old D3 = ROUND(C3/4,0) helper) is no longer needed and is removed
# - B3's cell format gains an extra decimal of precision (0.00 -> 0.000)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the raw ADC inputs; dependent formulas (C1, B3, C3) recalc automatically.
$ws.Range("A1").Value = 0.55000000000000004
$ws.Range("A3").Value = 10.65
$ws.Range("A4").Value = 1
$ws.Range("A6").Value = 1.1000000000000001

# Drop the now-unused column D (removes D3 and its formula, shrinks the used range).
$ws.Range("D1:D6").EntireColumn.Delete()

# Give the B3 ratio one more decimal of display precision.
$ws.Range("B3").NumberFormat = "0.000"
